$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.1170416182565993
$ws.Range("J2").Value = 0.1170416182565993
$ws.Range("M2").Value = 12.98850833333333
$ws.Range("N2").Value = 38.965525
$ws.Range("O2").Value = 0.4337249877124968
$ws.Range("P2").Value = 0.4337249877124967
$ws.Range("Q2").Value = 3.366171091711111
$ws.Range("R2").Value = 30.2955398254
$ws.Range("S2").Value = 0.05076387444019427
$ws.Range("T2").Value = 0.05076387444019427
$ws.Range("I3").Value = 0.1170416182565993
$ws.Range("J3").Value = 0.1170416182565993
$ws.Range("O3").Value = 0.4410094991449394
$ws.Range("P3").Value = 0.4410094991449393
$ws.Range("S3").Value = 0.05161646544645605
$ws.Range("T3").Value = 0.05161646544645605
$ws.Range("I4").Value = 0.1170416182565993
$ws.Range("J4").Value = 0.1170416182565993
$ws.Range("M4").Value = 0.5865036666666666
$ws.Range("N4").Value = 1.759511
$ws.Range("O4").Value = 0.01958510470101462
$ws.Range("P4").Value = 0.01958510470101462
$ws.Range("Q4").Value = 0.1520014182728889
$ws.Range("R4").Value = 1.368012764456
$ws.Range("S4").Value = 0.002292272347931682
$ws.Range("T4").Value = 0.002292272347931682
$ws.Range("I5").Value = 0.1170416182565993
$ws.Range("J5").Value = 0.1170416182565993
$ws.Range("M5").Value = 1.864071666666667
$ws.Range("N5").Value = 5.592214999999999
$ws.Range("O5").Value = 0.06224690626292447
$ws.Range("P5").Value = 0.06224690626292446
$ws.Range("Q5").Value = 0.4831027548488888
$ws.Range("R5").Value = 4.34792479364
$ws.Range("S5").Value = 0.007285478640479526
$ws.Range("T5").Value = 0.007285478640479526
$ws.Range("I6").Value = 0.1170416182565993
$ws.Range("J6").Value = 0.1170416182565993
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.300677666666667
$ws.Range("N6").Value = 3.902033
$ws.Range("O6").Value = 0.04343350217862475
$ws.Range("P6").Value = 0.04343350217862475
$ws.Range("Q6").Value = 0.3370905610408889
$ws.Range("R6").Value = 3.033815049368
$ws.Range("S6").Value = 0.005083527381537772
$ws.Range("T6").Value = 0.005083527381537772
$ws.Range("G7").Value = 1.708219666666666
$ws.Range("H7").Value = 5.124658999999999
$ws.Range("I7").Value = 0.7714488336573383
$ws.Range("J7").Value = 0.7714488336573383
$ws.Range("M7").Value = 12.98850833333333
$ws.Range("N7").Value = 38.965525
$ws.Range("O7").Value = 0.4337249877124968
$ws.Range("P7").Value = 0.4337249877124967
$ws.Range("Q7").Value = 22.18722537566389
$ws.Range("R7").Value = 199.685028380975
$ws.Range("S7").Value = 0.334596635898849
$ws.Range("T7").Value = 0.334596635898849
$ws.Range("G8").Value = 1.708219666666666
$ws.Range("H8").Value = 5.124658999999999
$ws.Range("I8").Value = 0.7714488336573383
$ws.Range("J8").Value = 0.7714488336573383
$ws.Range("O8").Value = 0.4410094991449394
$ws.Range("P8").Value = 0.4410094991449393
$ws.Range("Q8").Value = 22.55986495484889
$ws.Range("R8").Value = 203.03878459364
$ws.Range("S8").Value = 0.3402162637471705
$ws.Range("T8").Value = 0.3402162637471703
$ws.Range("G9").Value = 1.708219666666666
$ws.Range("H9").Value = 5.124658999999999
$ws.Range("I9").Value = 0.7714488336573383
$ws.Range("J9").Value = 0.7714488336573383
$ws.Range("M9").Value = 0.5865036666666666
$ws.Range("N9").Value = 1.759511
$ws.Range("O9").Value = 0.01958510470101462
$ws.Range("P9").Value = 0.01958510470101462
$ws.Range("Q9").Value = 1.001877097972111
$ws.Range("R9").Value = 9.016893881748999
$ws.Range("S9").Value = 0.01510890617865458
$ws.Range("T9").Value = 0.01510890617865458
$ws.Range("G10").Value = 1.708219666666666
$ws.Range("H10").Value = 5.124658999999999
$ws.Range("I10").Value = 0.7714488336573383
$ws.Range("J10").Value = 0.7714488336573383
$ws.Range("M10").Value = 1.864071666666667
$ws.Range("N10").Value = 5.592214999999999
$ws.Range("O10").Value = 0.06224690626292447
$ws.Range("P10").Value = 0.06224690626292446
$ws.Range("Q10").Value = 3.184243881076111
$ws.Range("R10").Value = 28.658194929685
$ws.Range("S10").Value = 0.04802030323531075
$ws.Range("T10").Value = 0.04802030323531074
$ws.Range("G11").Value = 1.708219666666666
$ws.Range("H11").Value = 5.124658999999999
$ws.Range("I11").Value = 0.7714488336573383
$ws.Range("J11").Value = 0.7714488336573383
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.300677666666667
$ws.Range("N11").Value = 3.902033
$ws.Range("O11").Value = 0.04343350217862475
$ws.Range("P11").Value = 0.04343350217862475
$ws.Range("Q11").Value = 2.221843170194111
$ws.Range("R11").Value = 19.996588531747
$ws.Range("S11").Value = 0.03350672459735353
$ws.Range("T11").Value = 0.03350672459735352
$ws.Range("G12").Value = 0.1229426666666667
$ws.Range("H12").Value = 0.368828
$ws.Range("I12").Value = 0.05552211970009493
$ws.Range("J12").Value = 0.05552211970009493
$ws.Range("M12").Value = 12.98850833333333
$ws.Range("N12").Value = 38.965525
$ws.Range("O12").Value = 0.4337249877124968
$ws.Range("P12").Value = 0.4337249877124967
$ws.Range("Q12").Value = 1.596841850522222
$ws.Range("R12").Value = 14.3715766547
$ws.Range("S12").Value = 0.02408133068469545
$ws.Range("T12").Value = 0.02408133068469544
$ws.Range("G13").Value = 0.1229426666666667
$ws.Range("H13").Value = 0.368828
$ws.Range("I13").Value = 0.05552211970009493
$ws.Range("J13").Value = 0.05552211970009493
$ws.Range("O13").Value = 0.4410094991449394
$ws.Range("P13").Value = 0.4410094991449393
$ws.Range("Q13").Value = 1.623661178542222
$ws.Range("R13").Value = 14.61295060688
$ws.Range("S13").Value = 0.02448578220040424
$ws.Range("T13").Value = 0.02448578220040423
$ws.Range("G14").Value = 0.1229426666666667
$ws.Range("H14").Value = 0.368828
$ws.Range("I14").Value = 0.05552211970009493
$ws.Range("J14").Value = 0.05552211970009493
$ws.Range("M14").Value = 0.5865036666666666
$ws.Range("N14").Value = 1.759511
$ws.Range("O14").Value = 0.01958510470101462
$ws.Range("P14").Value = 0.01958510470101462
$ws.Range("Q14").Value = 0.07210632478977777
$ws.Range("R14").Value = 0.6489569231079999
$ws.Range("S14").Value = 0.001087406527548626
$ws.Range("T14").Value = 0.001087406527548626
$ws.Range("G15").Value = 0.1229426666666667
$ws.Range("H15").Value = 0.368828
$ws.Range("I15").Value = 0.05552211970009493
$ws.Range("J15").Value = 0.05552211970009493
$ws.Range("M15").Value = 1.864071666666667
$ws.Range("N15").Value = 5.592214999999999
$ws.Range("O15").Value = 0.06224690626292447
$ws.Range("P15").Value = 0.06224690626292446
$ws.Range("Q15").Value = 0.2291739415577778
$ws.Range("R15").Value = 2.06256547402
$ws.Range("S15").Value = 0.003456080180490681
$ws.Range("T15").Value = 0.003456080180490681
$ws.Range("G16").Value = 0.1229426666666667
$ws.Range("H16").Value = 0.368828
$ws.Range("I16").Value = 0.05552211970009493
$ws.Range("J16").Value = 0.05552211970009493
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.300677666666667
$ws.Range("N16").Value = 3.902033
$ws.Range("O16").Value = 0.04343350217862475
$ws.Range("P16").Value = 0.04343350217862475
$ws.Range("Q16").Value = 0.1599087808137778
$ws.Range("R16").Value = 1.439179027324
$ws.Range("S16").Value = 0.002411520106955937
$ws.Range("T16").Value = 0.002411520106955937
$ws.Range("G17").Value = 0.123973
$ws.Range("H17").Value = 0.371919
$ws.Range("I17").Value = 0.05598742838596747
$ws.Range("J17").Value = 0.05598742838596747
$ws.Range("M17").Value = 12.98850833333333
$ws.Range("N17").Value = 38.965525
$ws.Range("O17").Value = 0.4337249877124968
$ws.Range("P17").Value = 0.4337249877124967
$ws.Range("Q17").Value = 1.610224343608333
$ws.Range("R17").Value = 14.492019092475
$ws.Range("S17").Value = 0.02428314668875803
$ws.Range("T17").Value = 0.02428314668875803
$ws.Range("G18").Value = 0.123973
$ws.Range("H18").Value = 0.371919
$ws.Range("I18").Value = 0.05598742838596747
$ws.Range("J18").Value = 0.05598742838596747
$ws.Range("O18").Value = 0.4410094991449394
$ws.Range("P18").Value = 0.4410094991449393
$ws.Range("Q18").Value = 1.637268433693333
$ws.Range("R18").Value = 14.73541590324
$ws.Range("S18").Value = 0.02469098775090868
$ws.Range("T18").Value = 0.02469098775090867
$ws.Range("G19").Value = 0.123973
$ws.Range("H19").Value = 0.371919
$ws.Range("I19").Value = 0.05598742838596747
$ws.Range("J19").Value = 0.05598742838596747
$ws.Range("M19").Value = 0.5865036666666666
$ws.Range("N19").Value = 1.759511
$ws.Range("O19").Value = 0.01958510470101462
$ws.Range("P19").Value = 0.01958510470101462
$ws.Range("Q19").Value = 0.07271061906766667
$ws.Range("R19").Value = 0.654395571609
$ws.Range("S19").Value = 0.001096519646879731
$ws.Range("T19").Value = 0.001096519646879731
$ws.Range("G20").Value = 0.123973
$ws.Range("H20").Value = 0.371919
$ws.Range("I20").Value = 0.05598742838596747
$ws.Range("J20").Value = 0.05598742838596747
$ws.Range("M20").Value = 1.864071666666667
$ws.Range("N20").Value = 5.592214999999999
$ws.Range("O20").Value = 0.06224690626292447
$ws.Range("P20").Value = 0.06224690626292446
$ws.Range("Q20").Value = 0.2310945567316667
$ws.Range("R20").Value = 2.079851010585
$ws.Range("S20").Value = 0.003485044206643514
$ws.Range("T20").Value = 0.003485044206643513
$ws.Range("G21").Value = 0.123973
$ws.Range("H21").Value = 0.371919
$ws.Range("I21").Value = 0.05598742838596747
$ws.Range("J21").Value = 0.05598742838596747
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1.300677666666667
$ws.Range("N21").Value = 3.902033
$ws.Range("O21").Value = 0.04343350217862475
$ws.Range("P21").Value = 0.04343350217862475
$ws.Range("Q21").Value = 0.1612489123696667
$ws.Range("R21").Value = 1.451240211327
$ws.Range("S21").Value = 0.002431730092777515
$ws.Range("T21").Value = 0.002431730092777515
